$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 values ---
$ws.Range("A2").Value = "MCH234-1"
$ws.Range("C2").Value = "SYMPOSIUM PROGRAMME, ART TOWARD SOCIAL DEVELOPMENT, CINEMA FOR LIBERATION, BRINGING THE STRUGGLE INTO FOCUS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1C | GRAP COUNT NUMER: NONE"

# --- Row 3 values ---
$ws.Range("A3").Value = "MCH234-2"
$ws.Range("C3").Value = "FESTIVAL EDITION, MEDU NEWSLETTER VOL.1 NO 2 & 4, A MEDU ART ENSAMBLE PRODUCTIAL"
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: CABINET 1C | GRAP COUNT NUMER: NONE"

# --- Apply formatting (font: Calibri 10, automatic/text1 color) to the two
# row ranges (column B intentionally excluded, matches target layout) ---
$rng1 = $ws.Range("A2:A3")
$rng1.Font.Name = "Calibri"
$rng1.Font.Size = 10
$rng1.Font.ThemeColor = 1

$rng2 = $ws.Range("C2:H3")
$rng2.Font.Name = "Calibri"
$rng2.Font.Size = 10
$rng2.Font.ThemeColor = 1
